# Loan RBI, Variable Instalments
#
# - On the "Repayment schedule" sheet, insert a new blank column in front of
#   column N. This pushes the existing "Outstanding" block (Due / Paid /
#   In Advance / Late) one column to the right (N->O, O->P, P->Q).
# - Make "Repayment schedule" the active sheet/tab and leave the selection on
#   R9, matching the workbook's stored view state.

$wb = $excel.ActiveWorkbook
$repay = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column N (14th column).
$repay.Columns.Item(14).Insert()

# Match the width carried onto the freshly inserted column (closest
# attainable value to the source column's 10.7109375 character width).
$repay.Columns.Item(14).ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and select R9, as last left
# by the author.
$repay.Activate()
$repay.Range("R9").Select() | Out-Null
